$wb = $excel.ActiveWorkbook

$treeSheets = @("Tree 1", "Tree 2", "Tree 3", "Tree 4", "Tree 5", "Tree 6", "Tree 7", "Tree 8")

foreach ($name in $treeSheets) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Rows.Item(1).Insert()
}
